# Solve Leetcode - 621. Task Scheduler - Dict, MaxHeap and DQ Approach
# Adds three new rows (13-15) to the "Neetcode 150" sheet under the Heap/PQ
# category: 973. K Closest Points to Origin, 215. Kth Largest Element in an
# Array, and 621. Task Scheduler.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Neetcode 150")

# ---------------------------------------------------------------------
# Row 13 - 973. K Closest Points to Origin
# ---------------------------------------------------------------------
$ws.Range("A13").Value = "Heap/PQ"
$ws.Range("B13").Value = "Medium"

$url13 = "https://leetcode.com/problems/k-closest-points-to-origin/"
$ws.Hyperlinks.Add($ws.Range("C13"), $url13, "", "", $url13) | Out-Null
$ws.Range("C13").Value = "973. K Closest Points to Origin"
$ws.Range("C13").Style = "Neutral"

$rsquo = [char]0x2019
$note1314 = "Quick - Use a heap to track top k elements (remove any elements after k for optmization). Store elements as negative since we need to find top k.`nOptimized - In Place Quick Select (Hoare's) (REVISE IT IF YOU DON" + $rsquo + "T REMEMBER CAUSE ITS TOO LONG TO DESCRIBE HERE, CHECK GITHUB REPO)"
$ws.Range("D13").Value = $note1314
$ws.Rows.Item(13).RowHeight = 57.6

# ---------------------------------------------------------------------
# Row 14 - 215. Kth Largest Element in an Array
# ---------------------------------------------------------------------
$ws.Range("A14").Value = "Heap/PQ"
$ws.Range("B14").Value = "Medium"

$url14 = "https://leetcode.com/problems/kth-largest-element-in-an-array/"
$ws.Hyperlinks.Add($ws.Range("C14"), $url14, "", "", $url14) | Out-Null
$ws.Range("C14").Value = "215. Kth Largest Element in an Array"
$ws.Range("C14").Style = "Neutral"

$ws.Range("D14").Value = $note1314
$ws.Rows.Item(14).RowHeight = 57.6

# ---------------------------------------------------------------------
# Row 15 - 621. Task Scheduler
# (description entered before the problem title, to match the original
# authoring order that shared-string indices were created in)
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "Heap/PQ"
$ws.Range("B15").Value = "Medium"

$note15 = "We first count the number of occurences of each letter with a dict (using inbuilt python func)`nThen we take only the count of values and heapify it, because we need to wait n time before doing the same task again so we start by the max occuring task so that we can fit as many small tasks in between as we need`nThen we go through the top task, decrement it to mark completion, and add it to a queue as [count, time we can do it again] at the end of the queue`nIf the current time is 1 and n is 3 and we're doing task `"A`" which has a count of 5, we'll reduce it to 4 and increase current time to 2, and the next time we can do it again is 2 + 3 = 5, so we'll add (4, 5) to the queue`nBecause of this, at each step, we'll also check if there is any task that shoud re enter the heap and pop the top of queue if true"
$ws.Range("D15").Value = $note15
$ws.Rows.Item(15).RowHeight = 115.2

$url15 = "https://leetcode.com/problems/task-scheduler/"
$ws.Hyperlinks.Add($ws.Range("C15"), $url15, "", "", $url15) | Out-Null
$ws.Range("C15").Value = "621. Task Scheduler"
$ws.Range("C15").Style = "Normal"

# ---------------------------------------------------------------------
# Update view: selection on the newly-added last row, scrolled into view
# ---------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("C15").Select() | Out-Null
